# aggiornamento a 9/09 compreso
# Add rows 367-374 (2021-09-02 through 2021-09-09) to Sheet1, continuing the
# existing daily covid-style series in columns A-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 367; A = 44441; B = 0; C = 2; D = 47.65308553728854 },
    @{ Row = 368; A = 44442; B = 1; C = 2; D = 47.65308553728854 },
    @{ Row = 369; A = 44443; B = 0; C = 2; D = 47.65308553728854 },
    @{ Row = 370; A = 44444; B = 0; C = 2; D = 47.65308553728854 },
    @{ Row = 371; A = 44445; B = 2; C = 3; D = 71.47962830593281 },
    @{ Row = 372; A = 44446; B = 0; C = 3; D = 71.47962830593281 },
    @{ Row = 373; A = 44447; B = 0; C = 3; D = 71.47962830593281 },
    @{ Row = 374; A = 44448; B = 0; C = 3; D = 71.47962830593281 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Copy formatting from the row above (carries over the date-formatted,
    # centered/bordered style used on column A) before writing new values.
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}

$excel.CutCopyMode = 0
